# Implemented timeseries to supim file, demand file next
#
# 1. SupIm sheet: extend the supim timeseries from 3 rows (hours 0-1) to
#    14 rows (hours 0-12), repeating the pattern of row 3 (hour 1) for the
#    new hours 2..12, and make it the active sheet/selection.
# 2. Process sheet: merge the two conditional-formatting rules covering
#    the new-process rows (11-13) into a single rule over A11:C13, and
#    drop tabSelected since SupIm becomes the active tab instead.
# 3. Process-Commodity sheet: conditional formatting left logically as-is
#    (its dxf slot shifts down automatically once the unused one is
#    freed up by the Process-sheet rule deletion above).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# SupIm: append rows 4..14 (hour index 2..12), copying the format + values
# of row 3 (hour index 1) which already hold the steady-state values.
# ---------------------------------------------------------------------
$wsSupIm = $wb.Worksheets.Item("SupIm")

for ($row = 4; $row -le 14; $row++) {
    $hour = $row - 2
    $wsSupIm.Range("A$row").Value = $hour
    $wsSupIm.Range("B$row").Value = 0.481
    $wsSupIm.Range("C$row").Value = 0.3
    $wsSupIm.Range("D$row").Value = 0.207

    # Copy row 3's formatting (keeps the A/D column styles) onto the new
    # row without disturbing the values just written above.
    $wsSupIm.Range("A3:D3").Copy()
    $wsSupIm.Range("A$row`:D$row").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Process sheet: combine the "A12:C13 A11 C11" + "B11" conditional
# formatting rules into the single "A11:C13" rule.
# ---------------------------------------------------------------------
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Range("B11").FormatConditions.Item(1).Delete()
$mergedRule = $wsProcess.Range("A12:C13").FormatConditions.Item(1)
$mergedRule.ModifyAppliesToRange($wsProcess.Range("A11:C13"))
$mergedRule.SetFirstPriority()

# ---------------------------------------------------------------------
# Make SupIm the active sheet/tab with the new selection, and restore
# Process' own selection (it is no longer the active tab).
# ---------------------------------------------------------------------
$wsProcess.Range("B11").Select()
$wsSupIm.Activate()
$wsSupIm.Range("K18").Select()

Write-Output "done"
